$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data (and its
# formatting) right by one column.
$ws.Columns.Item(1).Insert()

# The inserted column is blank/unformatted; copy the header formatting from
# the (now adjacent) former header cell so the new "id" header matches.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# Set the new header cell text.
$ws.Range("A1").Value = "id"

# Populate the new id column values for each data row.
$ws.Range("A2").Value = 8
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 10
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 1

# Append text to the description cell for row 4 (now column E after the shift).
$ws.Range("E4").Value = $ws.Range("E4").Value2 + " @asde"
